$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the duplicate "nip" (employee id) / "nama" (name) entry in row 5.
# Format column B as Text first so the numeric-looking id is stored as a
# shared string rather than being coerced into a number.
$ws.Range("B1:B1048576").NumberFormat = "@"
$ws.Range("B5").Value = "196512421107012111"
$ws.Range("C5").Value = "NARMINAL"

# Highlight duplicate values in the "nip" column (column B) so that this
# kind of mistake is caught before importing the batch data.
$colB = $ws.Range("B1:B1048576")
$dupRule = $colB.FormatConditions.AddUniqueValues()
$dupRule.DupeUnique = 1
$dupRule.Font.Color = 393372
$dupRule.Interior.Color = 13551615

# Move the active selection onto the cell that was just corrected.
$ws.Range("C5").Select()
